# add Use Item bug
# Append a new "ShowName" / "名字" row (row 20) to the Item table,
# matching the formatting of the preceding data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "ShowName"
$ws.Range("B20").Value = "string"
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = $false
$ws.Range("E20").Value = $false
$ws.Range("F20").Value = $true
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = "Friend"
$ws.Range("J20").Value = "名字"

# Match the text-formatted style (numFmtId "@") used by the other
# string-typed columns (A, B, I, J) in the existing rows.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("I20").NumberFormat = "@"
$ws.Range("J20").NumberFormat = "@"

# Move the active selection to the newly added cell, as in the target sheet.
[void]$ws.Range("A20").Select()
